$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contest 38 (row 50) - CSK vs KKR
$ws.Range("E50").Value = 60
$ws.Range("H50").Value = 0
$ws.Range("K50").Value = 80
$ws.Range("N50").Value = 40
$ws.Range("Q50").Value = 70
$ws.Range("T50").Value = 100
$ws.Range("W50").Value = 50
$ws.Range("Z50").Value = 20
$ws.Range("AC50").Value = 30

# Contest 39 (row 51) - RCB vs MI
$ws.Range("E51").Value = 70
$ws.Range("H51").Value = 80
$ws.Range("K51").Value = 60
$ws.Range("N51").Value = 20
$ws.Range("Q51").Value = 100
$ws.Range("T51").Value = 30
$ws.Range("W51").Value = 0
$ws.Range("Z51").Value = 50
$ws.Range("AC51").Value = 40

$excel.Calculate()
